$wb = $excel.ActiveWorkbook

# "Test Cases" sheet: D2 result changes from FAIL to PASS
$wsCases = $wb.Worksheets.Item("Test Cases")
$wsCases.Range("D2").Value = "PASS"

# "Test Steps" sheet: result column (H) updates for rows 20-28
$wsSteps = $wb.Worksheets.Item("Test Steps")

# H20 already holds a shared-string result value (FAIL) - just update it
$wsSteps.Range("H20").Value = "PASS"

# H21:H28 currently have no value (only a border style) - clear the
# existing formatting so the new cell matches the plain "PASS" cells,
# then set the value
$resultRange = $wsSteps.Range("H21:H28")
$resultRange.ClearFormats()
$resultRange.Value = "PASS"
